# Auto-generated Excel COM-interop script
# Applies the 2023-02-04 weekly crime data update across all affected sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 602
$ws.Cells.Item(3, 9).Value = 7487
$ws.Cells.Item(3, 10).Value = 658
$ws.Cells.Item(4, 5).Value = 1977
$ws.Cells.Item(4, 10).Value = 139
$ws.Cells.Item(5, 10).Value = 48
$ws.Cells.Item(6, 10).Value = 996
$ws.Cells.Item(7, 5).Value = 25980
$ws.Cells.Item(7, 9).Value = 26186
$ws.Cells.Item(7, 10).Value = 2443

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 27

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 27

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(7, 10).Value = 12

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 10).Value = 22
$ws.Cells.Item(6, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 91

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 10).Value = 38
$ws.Cells.Item(6, 10).Value = 31
$ws.Cells.Item(7, 10).Value = 91

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 31

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(6, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 50

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(6, 10).Value = 19
$ws.Cells.Item(7, 10).Value = 78
$ws.Cells.Item(8, 10).Value = 150
$ws.Cells.Item(9, 10).Value = 16
$ws.Cells.Item(11, 10).Value = 33
$ws.Cells.Item(12, 10).Value = 11
$ws.Cells.Item(14, 10).Value = 12
$ws.Cells.Item(15, 10).Value = 27
$ws.Cells.Item(18, 10).Value = 52
$ws.Cells.Item(19, 10).Value = 85
$ws.Cells.Item(20, 10).Value = 46
$ws.Cells.Item(23, 10).Value = 21
$ws.Cells.Item(25, 10).Value = 13
$ws.Cells.Item(29, 10).Value = 123
$ws.Cells.Item(32, 10).Value = 7
$ws.Cells.Item(33, 10).Value = 99
$ws.Cells.Item(34, 10).Value = 15
$ws.Cells.Item(37, 10).Value = 91
$ws.Cells.Item(42, 10).Value = 104
$ws.Cells.Item(43, 10).Value = 34
$ws.Cells.Item(44, 10).Value = 20
$ws.Cells.Item(45, 10).Value = 3
$ws.Cells.Item(51, 10).Value = 33
$ws.Cells.Item(52, 9).Value = 592
$ws.Cells.Item(52, 10).Value = 59
$ws.Cells.Item(53, 10).Value = 26
$ws.Cells.Item(54, 10).Value = 37
$ws.Cells.Item(55, 10).Value = 30
$ws.Cells.Item(60, 10).Value = 14
$ws.Cells.Item(63, 5).Value = 323
$ws.Cells.Item(65, 10).Value = 50
$ws.Cells.Item(66, 10).Value = 5
$ws.Cells.Item(67, 10).Value = 91
$ws.Cells.Item(73, 10).Value = 22
$ws.Cells.Item(74, 10).Value = 3
$ws.Cells.Item(76, 10).Value = 42
$ws.Cells.Item(79, 10).Value = 81
$ws.Cells.Item(83, 10).Value = 52
$ws.Cells.Item(84, 10).Value = 31
$ws.Cells.Item(85, 10).Value = 102
$ws.Cells.Item(87, 10).Value = 12
$ws.Cells.Item(89, 10).Value = 27
$ws.Cells.Item(95, 10).Value = 47
$ws.Cells.Item(96, 10).Value = 27
$ws.Cells.Item(100, 10).Value = 3
$ws.Cells.Item(101, 5).Value = 25980
$ws.Cells.Item(101, 9).Value = 26186
$ws.Cells.Item(101, 10).Value = 2443

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 10).Value = 13
$ws.Cells.Item(6, 10).Value = 23
$ws.Cells.Item(7, 10).Value = 52

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(6, 10).Value = 16
$ws.Cells.Item(7, 10).Value = 47

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 23
$ws.Cells.Item(3, 10).Value = 26
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(7, 10).Value = 99

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(6, 10).Value = 19
$ws.Cells.Item(7, 10).Value = 37

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 38
$ws.Cells.Item(7, 10).Value = 123

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 85

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 20

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 10).Value = 11
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(6, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 42

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 23
$ws.Cells.Item(3, 10).Value = 34
$ws.Cells.Item(7, 10).Value = 102

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 19

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 10).Value = 69
$ws.Cells.Item(7, 10).Value = 104

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(6, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 30

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(6, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 21

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 10).Value = 24
$ws.Cells.Item(3, 10).Value = 24
$ws.Cells.Item(7, 10).Value = 81

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(3, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 46

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(6, 10).Value = 37
$ws.Cells.Item(7, 10).Value = 52

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(6, 10).Value = 3

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(3, 9).Value = 193
$ws.Cells.Item(3, 10).Value = 20
$ws.Cells.Item(7, 9).Value = 592
$ws.Cells.Item(7, 10).Value = 59

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(6, 10).Value = 10
$ws.Cells.Item(7, 10).Value = 15

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 13

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 27

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(7, 10).Value = 5

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(3, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 33

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(6, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 16

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(7, 10).Value = 22

$ws = $wb.Worksheets.Item('Galewood')
$ws.Cells.Item(2, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 7

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 47
$ws.Cells.Item(3, 10).Value = 48
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(7, 10).Value = 150

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(3, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 33

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(6, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 14

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(2, 10).Value = 7
$ws.Cells.Item(6, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 34

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 10).Value = 17
$ws.Cells.Item(7, 10).Value = 26

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Cells.Item(6, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 3

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 18
$ws.Cells.Item(3, 10).Value = 33
$ws.Cells.Item(6, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 78

$ws = $wb.Worksheets.Item('Beverly')
$ws.Cells.Item(6, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 11

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(6, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 12

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(7, 10).Value = 3
